$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 397.1111
$ws.Range("J4").Value = 613
$ws.Range("L4").Value = 613
$ws.Range("N4").Value = -841
$ws.Range("H17").Value = 2770.3103
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H39").Value = 365.7
$ws.Range("I39").Value = 76.166664
$ws.Range("K39").Value = 228.499992
$ws.Range("M39").Value = 67.50000800000001
$ws.Range("H40").Value = 4081.125
$ws.Range("I40").Value = 3446.125
$ws.Range("J40").Value = 5351.125
$ws.Range("K40").Value = 3446.125
$ws.Range("L40").Value = 5351.125
$ws.Range("M40").Value = -3271.125
$ws.Range("N40").Value = -5701.125
$ws.Range("H42").Value = 180.14285
$ws.Range("I42").Value = 180.14285
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 540.4285500000001
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -310.4285500000001
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 700
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 2100
$ws.Range("M46").Value = -1981
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10972
$ws.Range("H60").Value = 700
$ws.Range("I60").Value = 700
$ws.Range("K60").Value = 2100
$ws.Range("M60").Value = -1616
$ws.Range("H104").Value = 1025.8
$ws.Range("I104").Value = 1025.8
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 3077.4
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -1330.4
$ws.Range("N104").ClearContents()
$ws.Range("H137").Value = 3077.1765
$ws.Range("I137").Value = 1718.3334
$ws.Range("K137").Value = 5155.0002
$ws.Range("M137").Value = -2605.0002
$ws.Range("H141").Value = 1408
$ws.Range("I141").Value = 1408
$ws.Range("K141").Value = 4224
$ws.Range("M141").Value = 956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3714.889
$ws.Range("I45").Value = 2063.3333
$ws.Range("J45").Value = 4540.6665
$ws.Range("K45").Value = 2063.3333
$ws.Range("L45").Value = 4540.6665
$ws.Range("M45").Value = -1686.3333
$ws.Range("N45").Value = -5294.6665
$ws.Range("H63").Value = 1142.1666
$ws.Range("J63").Value = 1124.5
$ws.Range("L63").Value = 1124.5
$ws.Range("N63").Value = -2496.5
$ws.Range("H66").Value = 1142.1666
$ws.Range("J66").Value = 1124.5
$ws.Range("L66").Value = 5622.5
$ws.Range("N66").Value = -12486.5
$ws.Range("H122").Value = 2582.3635
$ws.Range("I122").Value = 2388
$ws.Range("K122").Value = 7164
$ws.Range("M122").Value = -4714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H94").Value = 335.55554
$ws.Range("I94").Value = 277.5
$ws.Range("K94").Value = 277.5
$ws.Range("M94").Value = 173.5
$ws.Range("H99").Value = 2547.4
$ws.Range("I99").Value = 2681.5
$ws.Range("K99").Value = 2681.5
$ws.Range("M99").Value = -1183.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H31").Value = 8520.137000000001
$ws.Range("J31").Value = 9065.210999999999
$ws.Range("L31").Value = 9065.210999999999
$ws.Range("N31").Value = -9655.210999999999
$ws.Range("H32").Value = 987
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 8520.137000000001
$ws.Range("J34").Value = 9065.210999999999
$ws.Range("L34").Value = 9065.210999999999
$ws.Range("N34").Value = -9469.210999999999
$ws.Range("H50").Value = 20000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("M50").ClearContents()
$ws.Range("H51").Value = 44921.668
$ws.Range("J51").Value = 44921.668
$ws.Range("L51").Value = 44921.668
$ws.Range("N51").Value = -46393.668
$ws.Range("H60").Value = 12475
$ws.Range("I60").Value = 12475
$ws.Range("K60").Value = 12475
$ws.Range("M60").Value = -11964
$ws.Range("H61").Value = 44921.668
$ws.Range("J61").Value = 44921.668
$ws.Range("L61").Value = 44921.668
$ws.Range("N61").Value = -45617.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H87").Value = 5500
$ws.Range("I87").Value = 5500
$ws.Range("K87").Value = 16500
$ws.Range("M87").Value = -15252
$ws.Range("H90").Value = 5500
$ws.Range("I90").Value = 5500
$ws.Range("K90").Value = 49500
$ws.Range("M90").Value = -43260
$ws.Range("H101").Value = 8500
$ws.Range("I101").Value = 5000
$ws.Range("J101").Value = 12000
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 36000
$ws.Range("M101").Value = -12566
$ws.Range("N101").Value = -40868
$ws.Range("H102").Value = 14994.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 14994.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 44983.5
$ws.Range("N102").Value = -49851.5
$ws.Range("M102").ClearContents()
$ws.Range("H107").Value = 709.8570999999999
$ws.Range("I107").Value = 490
$ws.Range("J107").Value = 874.75
$ws.Range("K107").Value = 1470
$ws.Range("L107").Value = 2624.25
$ws.Range("M107").Value = 450
$ws.Range("N107").Value = -6464.25
$ws.Range("H118").Value = 300
$ws.Range("I118").Value = 300
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 900
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 343
$ws.Range("N118").ClearContents()
$ws.Range("H129").Value = 1487.375
$ws.Range("I129").Value = 583.9167
$ws.Range("K129").Value = 1751.7501
$ws.Range("M129").Value = 3248.2499
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 607.2857
$ws.Range("I16").Value = 375.16666
$ws.Range("K16").Value = 375.16666
$ws.Range("M16").Value = -205.16666
$ws.Range("H93").Value = 1977.4
$ws.Range("I93").Value = 1938
$ws.Range("J93").Value = 2069.3333
$ws.Range("K93").Value = 1938
$ws.Range("L93").Value = 2069.3333
$ws.Range("M93").Value = -690
$ws.Range("N93").Value = -4565.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 7400
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 7400
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 7400
$ws.Range("N33").Value = -7900
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 7400
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 7400
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 7400
$ws.Range("N36").Value = -7900
$ws.Range("M36").ClearContents()
$ws.Range("H52").Value = 46771
$ws.Range("I52").Value = 15042
$ws.Range("J52").Value = 78500
$ws.Range("K52").Value = 15042
$ws.Range("L52").Value = 78500
$ws.Range("M52").Value = -14816
$ws.Range("N52").Value = -78952
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
